$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

try {
  $win = $wb.Windows.Item(1)
  Write-Host "win: $win"
  Write-Host "left: $($win.Left)"
  $win.Left = 780
  $win.Top = 780
  Write-Host "set ok"
} catch {
  Write-Host "err: $_"
}
